# Locate the list paragraph that ends with "Crear Controller para cada
# clase. (BackEnd)" - the new task item must be inserted right after it.
$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Crear Controller para cada clase*(BackEnd)*") {
        $target = $p
    }
}

# Insert a brand-new empty paragraph right after the target paragraph.
# It inherits the pPr (pStyle "Prrafodelista", numPr ilvl 0 / numId 5) of
# the target paragraph, matching the list formatting used throughout the
# "Tareas" list.
$target.Range.InsertParagraphAfter()

# The newly created paragraph occupies exactly the single pilcrow
# character right after the end of the target paragraph's range.
$newParaStart = $target.Range.End
$newParaEnd = $newParaStart + 1
$insPoint = $d.Range($newParaStart, $newParaEnd)

# Build the run-level content for the new paragraph, replicating the
# exact OOXML runs (including the spell-check proofErr markers that wrap
# the English word "login"), and insert it as raw WordprocessingML.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Implementar JWT para el sistema de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>login</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>'

$insPoint.InsertXML($xml)
